$d = $word.ActiveDocument

# --- Paragraph 1: add paragraph border (space=5 on all sides) ---
$p1 = $d.Paragraphs.Item(1)
$p1.Format.Borders.DistanceFromTop = 5
$p1.Format.Borders.DistanceFromLeft = 5
$p1.Format.Borders.DistanceFromBottom = 5
$p1.Format.Borders.DistanceFromRight = 5

# --- Paragraph 1: change left indent 120 -> 225 twips (6pt -> 11.25pt) ---
$p1.Format.LeftIndent = 225 / 20.0

# --- Paragraph 1: replace placeholder text and drop the trailing-space run ---
$rng = $p1.Range
$rng.Find.Execute("**ID__AFFARS_5319_topic_11__ID** ", $true, $false, $false, $false, $false,
                   $true, 1, $false, "**ID__AFFARS_SUBPART_5319_13__ID**", 2)

Write-Host "done"
